$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for the "family with disabilities Persons" data,
#     right after the years row (row 3), before the old "Number of disability
#     persons" row (row 4). This shifts the old rows 4 and 5 down to 5 and 6.
$ws.Rows.Item(4).Insert()

# --- Row 1: title (merged A1:I1) ---
$ws.Range("A1:I1").Merge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Bolnisi Municipality"
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# --- Row 2: "(End of year, persons)" ---
$ws.Range("A2").Value = "(End of year, persons)"
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.ColorIndex = 1
$ws.Range("A2").Interior.Pattern = 1
$ws.Range("A2").Interior.ThemeColor = 1
$ws.Rows.Item(2).RowHeight = 14.5

# --- Row 3: years header, A3 style only (blank label cell) ---
$ws.Range("A3").Font.Name = "Arial"
$ws.Range("A3").Font.Size = 10
$ws.Range("A3").Font.ThemeColor = 1

# --- Row 4: new "family with disabilities Persons" data row ---
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 9
$ws.Range("A4").Font.ColorIndex = 1
$ws.Range("A4").Interior.Pattern = 1
$ws.Range("A4").Interior.ThemeColor = 1
$ws.Range("A4").Borders.Item(9).LineStyle = 1
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 24.75

$vals4 = @(471, 448, 489, 545, 605, 673, 725, 762)
$cols = @("B", "C", "D", "E", "F", "G", "H", "I")
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "4")
    $cell.Value = $vals4[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.ColorIndex = 1
    $cell.Interior.Pattern = 1
    $cell.Interior.ThemeColor = 1
}

# --- Row 5: "disabilities Persons" data row (was "Number of disability persons") ---
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 9
$ws.Range("A5").Font.ColorIndex = 1
$ws.Range("A5").Interior.Pattern = 1
$ws.Range("A5").Interior.ThemeColor = 1
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 21

$vals5 = @(522, 498, 536, 602, 664, 735, 792, 824)
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "5")
    $cell.Value = $vals5[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.ColorIndex = 1
    $cell.Interior.Pattern = 1
    $cell.Interior.ThemeColor = 1
}
# I5 additionally carries a bottom border under the new layout
$ws.Range("I5").Borders.Item(9).LineStyle = 1

# --- Row 6: Source row (was row 5), now merged A6:H6 ---
$ws.Range("A6:H6").Merge()
$ws.Range("A6").Font.Name = "Arial"
$ws.Range("A6").Font.Size = 9
$ws.Range("A6").Font.ColorIndex = 1
$ws.Range("A6").Interior.Pattern = 1
$ws.Range("A6").Interior.ThemeColor = 1
$ws.Range("A6").HorizontalAlignment = -4131
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("A6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 27.75

# --- Column A width update ---
$ws.Columns.Item(1).ColumnWidth = 20.8

# --- Selection on the title row ---
$ws.Range("A1:I1").Select()

Write-Host "edit applied"
